$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the "Label" header in column H, matching the style of the other headers (s="1")
$ws.Range("H1").Value = "Label"
$ws.Range("H1").Style = $ws.Range("G1").Style

# Add the Label column values (0 = Control, 1 = MDD) for both Iterations=100 and Iterations=200 blocks
$labels = @(0,0,0,0,0,1,1,1,1,1,0,0,0,0,0,1,1,1,1,1)
for ($i = 0; $i -lt $labels.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 8).Value = $labels[$i]
}

# Refit the Prediction (D) / Error (E) columns with updated values from the re-run NCDE fit
$ws.Range("D2").Value = 0.5073717949588463
$ws.Range("E2").Value = 0.5073717949588463

$ws.Range("D3").Value = 0.5330848053201347
$ws.Range("E3").Value = 0.5330848053201347

$ws.Range("D4").Value = 0.5114138384180313
$ws.Range("E4").Value = 0.5114138384180313

$ws.Range("D5").Value = 0.4567653124054399
$ws.Range("E5").Value = 0.4567653124054399

$ws.Range("D6").Value = 0.2591942689264033
$ws.Range("E6").Value = 0.2591942689264033

$ws.Range("D7").Value = 0.5226258795120079
$ws.Range("E7").Value = 0.4773741204879921

$ws.Range("D8").Value = 0.4872523005617599
$ws.Range("E8").Value = 0.5127476994382401

$ws.Range("D9").Value = 0.5142405742298127
$ws.Range("E9").Value = 0.4857594257701873

$ws.Range("D10").Value = 0.5406530646132104
$ws.Range("E10").Value = 0.4593469353867896

$ws.Range("D11").Value = 0.5292507757369816
$ws.Range("E11").Value = 0.4707492242630184
